$wb = $excel.ActiveWorkbook

# --- ManageOrdersPage: drop the "Time Field"/"Expected Text" test columns
#     that were only there by error, and renumber the Order Id column ---
$orders = $wb.Worksheets.Item("ManageOrdersPage")

# Row 1 headers: clear B1/C1 text (keep their bold style), drop D1 entirely.
$orders.Range("B1").Value = ""
$orders.Range("C1").Value = ""
$orders.Range("D1").ClearContents()

# Row 2: renumber the Order Id, clear the time-field/expected-text values
# (D2 keeps its wrap-text style, just with no content).
$orders.Range("A2").Value = 136
$orders.Range("B2").Value = ""
$orders.Range("C2").Value = ""
$orders.Range("D2").Value = ""

# Rows 3-9: Order Id shifts down by 9 for each row.
$orders.Range("A3").Value = 135
$orders.Range("A4").Value = 134
$orders.Range("A5").Value = 133
$orders.Range("A6").Value = 132
$orders.Range("A7").Value = 131
$orders.Range("A8").Value = 130
$orders.Range("A9").Value = 129

# ManageOrdersPage becomes the active/selected sheet & tab (was
# ManageOfferCodePage).
$orders.Activate()
